# Daily attendance processing - 2025-12-06 04:29:35
# Normalizes the "Recorded By" (column G) entries so that when the list of
# recorders contains "dnasr281@gmail.com" together with exactly one other
# recorder (e.g. "System" or "admin@admin.com"), the two names are swapped
# in the displayed, comma-separated order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "dnasr281@gmail.com,*") {
        $parts = $val -split ", "
        if ($parts.Count -eq 2) {
            $newVal = $parts[1] + ", " + $parts[0]
            $cell.Value2 = $newVal
        }
    }
}
